# Refresh scraped market-board derived values (currentAveragePrice*, Leve*Profit*)
# across the per-job Leve-profit sheets, as produced by the scheduled data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3940.818
$ws.Range("I62").Value = 3260.3333
$ws.Range("J62").Value = 7003
$ws.Range("K62").Value = 3260.3333
$ws.Range("L62").Value = 7003
$ws.Range("M62").Value = -2636.3333
$ws.Range("N62").Value = -8251

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3940.818
$ws.Range("I65").Value = 3260.3333
$ws.Range("J65").Value = 7003
$ws.Range("K65").Value = 16301.6665
$ws.Range("L65").Value = 35015
$ws.Range("M65").Value = -13181.6665
$ws.Range("N65").Value = -41255

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2727.25
$ws.Range("J88").Value = 2636.3333
$ws.Range("L88").Value = 2636.3333
$ws.Range("N88").Value = -3448.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2727.25
$ws.Range("J91").Value = 2636.3333
$ws.Range("L91").Value = 2636.3333
$ws.Range("N91").Value = -5444.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 983
$ws.Range("I94").Value = 983
$ws.Range("K94").Value = 983
$ws.Range("M94").Value = -532

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1577
$ws.Range("I111").Value = 721.5
$ws.Range("K111").Value = 2164.5
$ws.Range("M111").Value = 902.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1089649.2
$ws.Range("I112").Value = 1534.75
$ws.Range("K112").Value = 4604.25
$ws.Range("M112").Value = -3496.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3287
$ws.Range("I116").Value = 3004.5
$ws.Range("K116").Value = 3004.5
$ws.Range("M116").Value = 437.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 23811294
$ws.Range("I118").Value = 27778176
$ws.Range("K118").Value = 83334528
$ws.Range("M118").Value = -83332871

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2105.9443
$ws.Range("I132").Value = 2105.9443
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6317.8329
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3787.8329
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 585.8461
$ws.Range("I4").Value = 401.45456
$ws.Range("K4").Value = 401.45456
$ws.Range("M4").Value = -285.45456

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2276917.8
$ws.Range("I32").Value = 2329848.8
$ws.Range("K32").Value = 2329848.8
$ws.Range("M32").Value = -2329561.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2929.1428
$ws.Range("I110").Value = 1574.4706
$ws.Range("K110").Value = 1574.4706
$ws.Range("M110").Value = 470.5293999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 14707693
$ws.Range("I132").Value = 1298.3793
$ws.Range("K132").Value = 3895.1379
$ws.Range("M132").Value = -1365.1379

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2536.75
$ws.Range("I107").Value = 1778.75
$ws.Range("K107").Value = 1778.75
$ws.Range("M107").Value = 141.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4183.091
$ws.Range("I122").Value = 4290
$ws.Range("J122").Value = 3114
$ws.Range("K122").Value = 12870
$ws.Range("L122").Value = 9342
$ws.Range("M122").Value = -10420
$ws.Range("N122").Value = -14242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 699.5
$ws.Range("I75").Value = 398
$ws.Range("J75").Value = 800
$ws.Range("K75").Value = 1194
$ws.Range("L75").Value = 2400
$ws.Range("M75").Value = -196
$ws.Range("N75").Value = -4396

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 699.5
$ws.Range("I78").Value = 398
$ws.Range("J78").Value = 800
$ws.Range("K78").Value = 3582
$ws.Range("L78").Value = 7200
$ws.Range("M78").Value = 1410
$ws.Range("N78").Value = -17184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 19861
$ws.Range("I87").Value = 19861
$ws.Range("K87").Value = 59583
$ws.Range("M87").Value = -58335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 19861
$ws.Range("I90").Value = 19861
$ws.Range("K90").Value = 178749
$ws.Range("M90").Value = -172509

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 261308.14
$ws.Range("J131").Value = 281332.53
$ws.Range("L131").Value = 843997.5900000001
$ws.Range("N131").Value = -854077.5900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 9981.799999999999
$ws.Range("J15").Value = 9981.799999999999
$ws.Range("L15").Value = 9981.799999999999
$ws.Range("N15").Value = -10557.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 9981.799999999999
$ws.Range("J81").Value = 9981.799999999999
$ws.Range("L81").Value = 9981.799999999999
$ws.Range("N81").Value = -11977.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 9981.799999999999
$ws.Range("J84").Value = 9981.799999999999
$ws.Range("L84").Value = 29945.4
$ws.Range("N84").Value = -39929.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3666.6538
$ws.Range("I40").Value = 1663.4286
$ws.Range("J40").Value = 6003.75
$ws.Range("K40").Value = 1663.4286
$ws.Range("L40").Value = 6003.75
$ws.Range("M40").Value = -1527.4286
$ws.Range("N40").Value = -6275.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4125.6113
$ws.Range("J46").Value = 4882.7856
$ws.Range("L46").Value = 4882.7856
$ws.Range("N46").Value = -5258.7856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3958.1333
$ws.Range("I82").Value = 3697.4167
$ws.Range("J82").Value = 5001
$ws.Range("K82").Value = 3697.4167
$ws.Range("L82").Value = 5001
$ws.Range("M82").Value = -3336.4167
$ws.Range("N82").Value = -5723

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3958.1333
$ws.Range("I85").Value = 3697.4167
$ws.Range("J85").Value = 5001
$ws.Range("K85").Value = 3697.4167
$ws.Range("L85").Value = 5001
$ws.Range("M85").Value = -2449.4167
$ws.Range("N85").Value = -7497

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 16639.148
$ws.Range("I93").Value = 1187.1875
$ws.Range("K93").Value = 1187.1875
$ws.Range("M93").Value = 60.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 49498.5
$ws.Range("J96").Value = 49999
$ws.Range("L96").Value = 49999
$ws.Range("N96").Value = -55491

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3490.3845
$ws.Range("I122").Value = 3498.7144
$ws.Range("J122").Value = 3480.6667
$ws.Range("K122").Value = 10496.1432
$ws.Range("L122").Value = 10442.0001
$ws.Range("M122").Value = -8046.143199999999
$ws.Range("N122").Value = -15342.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4303.8184
$ws.Range("I132").Value = 2494.4119
$ws.Range("J132").Value = 10455.8
$ws.Range("K132").Value = 7483.2357
$ws.Range("L132").Value = 31367.4
$ws.Range("M132").Value = -4953.2357
$ws.Range("N132").Value = -36427.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2984.5454
$ws.Range("I81").Value = 2571.5
$ws.Range("K81").Value = 5143
$ws.Range("M81").Value = -4082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2984.5454
$ws.Range("I84").Value = 2571.5
$ws.Range("K84").Value = 25715
$ws.Range("M84").Value = -20411
